$d = $word.ActiveDocument
$t = $d.Tables(1)

# "Sprint No." value cell: row 2, column 4 -> change "1" to "2"
$sprintCell = $t.Cell(2, 4)
$sprintCell.Range.Text = "2"

# "Review Date" value cell: row 3, column 2 (spans the row) -> change "02/09/18" to "02/21/18"
$dateCell = $t.Cell(3, 2)
$dateCell.Range.Text = "02/21/18"
